# Atualizando o arquivo XLSX
# Apply updated odds values to Sheet1 as described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 1.6
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.95

# Row 3
$ws.Range("G3").Value = 1.8
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 2.3
$ws.Range("Q3").Value = 2.88
$ws.Range("R3").Value = 1.4

# Row 4
$ws.Range("G4").Value = 1.75
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 5.75
$ws.Range("J4").Value = 2.5
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 6
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("U4").Value = 2.38
$ws.Range("V4").Value = 1.53
$ws.Range("W4").Value = 5
$ws.Range("X4").Value = 7
$ws.Range("Y4").Value = 9.5
$ws.Range("AA4").Value = 19
$ws.Range("AC4").Value = 6
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 81
$ws.Range("AH4").Value = 11
$ws.Range("AJ4").Value = 19
$ws.Range("AM4").Value = 67
$ws.Range("AN4").Value = 3.5
$ws.Range("AP4").Value = 26
$ws.Range("AQ4").Value = 34
$ws.Range("AR4").Value = 67
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.25
$ws.Range("AU4").Value = 10
$ws.Range("AV4").Value = 81
$ws.Range("AW4").Value = 6.5
$ws.Range("AY4").Value = 41
$ws.Range("AZ4").Value = 126
$ws.Range("BA4").Value = 201

# Row 10
$ws.Range("N10").Value = 8
$ws.Range("BD10").Value = 126

# Row 13
$ws.Range("G13").Value = 1.25
$ws.Range("H13").Value = 5.25
$ws.Range("I13").Value = 13
$ws.Range("J13").Value = 1.67
$ws.Range("L13").Value = 9.5
$ws.Range("N13").Value = 12
$ws.Range("O13").Value = 1.2
$ws.Range("P13").Value = 4.33
$ws.Range("Q13").Value = 1.65
$ws.Range("R13").Value = 2.2
$ws.Range("U13").Value = 2.2
$ws.Range("V13").Value = 1.62
$ws.Range("X13").Value = 6
$ws.Range("AA13").Value = 12
$ws.Range("AB13").Value = 34
$ws.Range("AC13").Value = 12
$ws.Range("AE13").Value = 23
$ws.Range("AF13").Value = 81
$ws.Range("AH13").Value = 23
$ws.Range("AK13").Value = 151
$ws.Range("AL13").Value = 81
$ws.Range("AN13").Value = 3.2
$ws.Range("AO13").Value = 5.5
$ws.Range("AP13").Value = 19
$ws.Range("AS13").Value = 151
$ws.Range("AW13").Value = 11
$ws.Range("AX13").Value = 51
$ws.Range("AY13").Value = 51
$ws.Range("AZ13").Value = 251
$ws.Range("BA13").Value = 251
